# fix error tahap 3
# Update "Hasil Stemming" column (C) values per corrected stemming diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C10" = "kadang"
    "C11" = "dip"
    "C12" = "dip"
    "C13" = "erang"
    "C14" = "rtas"
    "C16" = "pajang"
    "C17" = "colok"
    "C24" = "pandang"
    "C26" = "arik"
    "C35" = "temu"
    "C37" = "nyelid"
    "C38" = "putus"
    "C39" = "bentuk"
    "C42" = "putus"
    "C43" = "ambah"
    "C50" = "letak"
    "C60" = "it"
    "C61" = "reta"
    "C62" = "jad"
    "C63" = "kenal"
    "C72" = "anja"
    "C74" = "tinggal"
    "C77" = "banding"
    "C81" = "pilik"
    "C85" = "lar"
    "C89" = "liar"
    "C93" = "jahat"
    "C96" = "kenal"
    "C103" = "keras"
    "C105" = "rkosa"
    "C106" = "culik"
    "C107" = "bunuh"
    "C111" = "liput"
    "C113" = "cekam"
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
